# Add 2022-Q3 data
#
# The workbook tracks one "总计" (summary) sheet plus one sheet per
# quarter. This change adds a new quarter (2022-Q3):
#   1. "总计" gains a new top data row for 2022-Q3, and the previously
#      oldest quarter (2020-Q4) gets appended as a new trailing row.
#   2. A brand-new "2022-Q3" worksheet is inserted right after "总计",
#      holding the per-fund detail rows for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift quarter labels/values up by one row and
#    insert the new 2022-Q3 figures at the top; the row that falls off
#    the bottom (2020-Q4) is appended as a brand new last row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$quarters = @(
    @{ B = "2022-Q3"; D = 0.12 },
    @{ B = "2022-Q2"; D = 0.14 },
    @{ B = "2022-Q1"; D = 0.19 },
    @{ B = "2021-Q4"; D = 0.24 },
    @{ B = "2021-Q3"; D = 0.21 },
    @{ B = "2021-Q2"; D = 0.21 },
    @{ B = "2021-Q1"; D = 0.18 },
    @{ B = "2020-Q4"; D = 0.09 }
)

# Bring the formatting of the newly-appended row 9 in line with the
# rest of the table (column A carries a centered/bordered style).
$total.Range("A8").Copy($total.Range("A9"))

for ($i = 0; $i -lt $quarters.Count; $i++) {
    $row = 2 + $i
    $total.Cells.Item($row, 1).Value = $i
    $total.Cells.Item($row, 2).Value = $quarters[$i].B
    $total.Cells.Item($row, 3).Value = 2
    $total.Cells.Item($row, 4).Value = $quarters[$i].D
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail sheet right after "总计", cloned
#    from the "2022-Q2" sheet so it inherits identical formatting, then
#    overwrite its figures with the 2022-Q3 numbers.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $total)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# D:G are stored as plain text (keeps the "78.58"-style fixed
# decimals), H is a genuine number - matches the other quarter sheets.
$figures = $q3.Range("D2:G3")
$figures.NumberFormat = "@"

$q3.Range("D2").Value = "1.58"
$q3.Range("E2").Value = "78.58"
$q3.Range("F2").Value = "3.79"
$q3.Range("G2").Value = "0.0599"
$q3.Range("H2").Value = 6

$q3.Range("D3").Value = "1.58"
$q3.Range("E3").Value = "78.58"
$q3.Range("F3").Value = "3.79"
$q3.Range("G3").Value = "0.0599"
$q3.Range("H3").Value = 6

# ---------------------------------------------------------------------
# Restore the originally-active sheet (last quarter, "2020-Q4" - still
# the final tab after the insertion above) as the selected one.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
